$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Featured Sessions")

# Insert a new row at position 4, shifting existing rows down
$ws.Rows.Item(4).Insert()

# Fill in the new row's content
$ws.Range("A4").Value = "Statistics in Practice: Simulation studies as a tool to evaluate and compare the properties of statistical methods – an overview"
$ws.Range("B4").Value = "Willi Sauerbrei"

# Apply wrap text alignment to the new title cell
$ws.Range("A4").WrapText = $true

# Update selection to match final state
$ws.Range("B6").Select()
